$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 9135
$ws1.Range("F6").Value = 9135
$ws1.Range("F10").Value = 219
$ws1.Range("F12").Value = 385
$ws1.Range("F16").Value = 11684
$ws1.Range("F17").Value = 11684
$ws1.Range("F18").Value = 51
$ws1.Range("F26").Value = 18
$ws1.Range("F37").Value = 4165
$ws1.Range("F39").Value = 3587
$ws1.Range("F40").Value = 323
$ws1.Range("F43").Value = 1294
$ws1.Range("F44").Value = 183
$ws1.Range("F46").Value = 387
$ws1.Range("F47").Value = 447
$ws1.Range("F50").Value = 115

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 10
$ws2.Range("F13").Value = 44
$ws2.Range("F14").Value = 21

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F10").Value = 9135
$ws4.Range("F11").Value = 9135
$ws4.Range("F15").Value = 219
$ws4.Range("F16").Value = 385
$ws4.Range("F19").Value = 11684
$ws4.Range("F20").Value = 11684
$ws4.Range("F23").Value = 10
$ws4.Range("F27").Value = 18
$ws4.Range("F42").Value = 3587
$ws4.Range("F45").Value = 1294
$ws4.Range("F46").Value = 183
$ws4.Range("F47").Value = 387
$ws4.Range("F49").Value = 447
